$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# 1. Move the existing "Simulation 1/2/3" parameter block (old columns C:E,
#    plus the demography legend that spilled into F:J on row 3) to the right,
#    into columns K:M (and N:R for the row-3 legend), to make room for a new
#    set of descriptive "range" columns C:G.
# ---------------------------------------------------------------------------

# Row 3 carries the demography legend across 8 columns (C:J) in the old
# layout; move the whole thing to K:R in one shot so formatting travels too.
$ws.Range("C3:J3").Copy($ws.Range("K3"))

# Every other row just has the 3 simulation columns C:E -> K:M.
$ws.Range("C1:E2").Copy($ws.Range("K1"))
$ws.Range("C4:E27").Copy($ws.Range("K4"))

# Wipe out the old D:J content now that it has been relocated (row 1's "Simulation 2"
# label in D1 is left in place -- it was not part of the relocated block).
$ws.Range("D2:J27").ClearContents()
$ws.Range("E1:J1").ClearContents()

# ---------------------------------------------------------------------------
# 2. Populate the new "Search range" column (C) with the updated values and
#    add the new demography-pattern legend in D3:G3.
# ---------------------------------------------------------------------------

$ws.Range("C5").Value = "1000:10000"
$ws.Range("C6").Value = "10^UNIF(1/N, 0)"
$ws.Range("C7").Value = "UNIF(0,1)"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = "10^runif(K, min = -8.5, max = -7.5)"
$ws.Range("C10").Value = "10^runif(K, min = -9, max = -7)"
$ws.Range("C11").Value = "round(10^runif(K, min = 0, max = 4))"

$ws.Range("B3").Value = "pattern of how N changed over generations"
$ws.Range("D3").Value = "growth"
$ws.Range("E3").Value = "decay"
$ws.Range("F3").Value = "cycling"
$ws.Range("G3").Value = "chaotic"

# f0 range differs only in the 3rd simulation column.
$ws.Range("M13").Value = "UNIF(1/(2N), 0.1)"

# tau: simulation 1 stays a constant 1, simulations 2/3 use the 10^UNIF draw.
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = "10^UNIF(0,3)"
$ws.Range("M11").Value = "10^UNIF(0,3)"

# sigma (selfing rate) search range is fixed at 0; simulation columns are 1.
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1

# ---------------------------------------------------------------------------
# 3. Cosmetic bits: row heights, used range, active selection.
# ---------------------------------------------------------------------------

$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 30
$ws.Rows.Item(21).RowHeight = 30
$ws.Rows.Item(22).RowHeight = 30
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30
$ws.Rows.Item(25).RowHeight = 30
$ws.Rows.Item(26).RowHeight = 30
$ws.Rows.Item(27).RowHeight = 45

$ws.Range("E4").Select()
